$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: column A stays the same, column B widens ---
$ws.Columns("B").ColumnWidth = 76.88671875

# --- Fill in the new cell values (order chosen to mirror the authoring session) ---
$ws.Range("D2").Value = "Low"
$ws.Range("D3").Value = "Medium"
$ws.Range("A3").Value = "Unify monthly sunshine data"
$ws.Range("D4").Value = "High"
$ws.Range("C4").Value = "No"
$ws.Range("A4").Value = "When cm is written in the days count we also multiply it by 10"
$ws.Range("B4").Value = "Check Zurich as an example. I tshould  happen in rainy and snowy days"
$ws.Range("B3").Value = "The records 'Mean daily sunshine hours', 'Percent possible sunshine' are displayed as 'Mean monthly sunshine hours'`nCheck [Zurich, Zunyi] as has the 2 fields so we can check it`nUse Freiburg it to compare the mean daily and mean monthly sunshine hours"
$ws.Range("A5").Value = "Get data from the"
$ws.Range("B6").Value = "https://ghsl.jrc.ec.europa.eu/datasets.php"

$ws.Range("C3").Value = "NO"

# --- Row heights for the new wrapped rows ---
$ws.Rows(3).RowHeight = 60.6
$ws.Rows(4).RowHeight = 28.8

# --- Wrap text, applied cell-by-cell so no empty placeholder cells are created ---
$ws.Range("A1").WrapText = $true
$ws.Range("B1").WrapText = $true
$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Range("A3").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("A4").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("B6").WrapText = $true

# --- Data validation dropdown for the Importance column ---
$range = $ws.Range("D2:D34")
$range.Validation.Add(3, 1, 1, '"High,Medium,Low"')
$range.Validation.IgnoreBlank = $true
$range.Validation.InCellDropdown = $true
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true

Write-Host "edit applied"
